$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 251, shifting rows 251:281 down to 252:282
$ws.Rows.Item(251).Insert()

# Populate the new row 251 with its values
$ws.Cells.Item(251, 1).Value = 10
$ws.Cells.Item(251, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(251, 3).Value = "La Araucanía"
$ws.Cells.Item(251, 4).Value = 44449
$ws.Cells.Item(251, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(251, 5).Value = 9
$ws.Cells.Item(251, 6).Value = 100112043
$ws.Cells.Item(251, 7).Value = "Pepino ensalada"
$ws.Cells.Item(251, 8).Value = "Sin especificar"
$ws.Cells.Item(251, 9).Value = "Primera"
$ws.Cells.Item(251, 10).Value = 195
$ws.Cells.Item(251, 11).Value = 17000
$ws.Cells.Item(251, 12).Value = 17000
$ws.Cells.Item(251, 13).Value = 17000
$ws.Cells.Item(251, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(251, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(251, 16).Value = 283
$ws.Cells.Item(251, 17).Value = 60
$ws.Cells.Item(251, 18).Value = "Hortaliza"
